# Daily attendance processing - 2026-02-07 04:14:59 UTC
# Updates the "Recorded By" text in column G from
#   "Miss Dina Nasr, Administrator"
# to
#   "Administrator, Miss Dina Nasr"
# for the specific rows that were touched by the upstream edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rows = @(4,5,8,17,18,19,20,21,22,23,24,26,30,31,34,43,44,45,46,47,48,49,50,52,56,57,60,69,70,71,72,73,74,75,76,78,80,81,82,93,94,96,99,101,106,107,108,119,120,122,125,127,132,133,134,145,146,148,151,153)

$oldText = "Miss Dina Nasr, Administrator"
$newText = "Administrator, Miss Dina Nasr"

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    if ($cell.Value2 -eq $oldText) {
        $cell.Value2 = $newText
    }
}
